$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '62.903.95'
Set-TextValue $ws.Range("E2") '  -0.62%  '
Set-TextValue $ws.Range("D3") '2.571.24'
Set-TextValue $ws.Range("E3") '  -0.20%  '
Set-TextValue $ws.Range("E4") '  +0.05%  '
Set-TextValue $ws.Range("D5") '583.21'
Set-TextValue $ws.Range("E5") '  -0.34%  '
Set-TextValue $ws.Range("D6") '144.05'
Set-TextValue $ws.Range("E6") '  -2.68%  '
Set-TextValue $ws.Range("E7") '  +0.05%  '
Set-TextValue $ws.Range("D8") '0.589'
Set-TextValue $ws.Range("E8") '  -2.15%  '
Set-TextValue $ws.Range("E9") '  -2.64%  '
Set-TextValue $ws.Range("E10") '  -1.19%  '
Set-TextValue $ws.Range("E11") '  -0.44%  '
Set-TextValue $ws.Range("E12") '  -2.48%  '
Set-TextValue $ws.Range("D13") '27.00'
Set-TextValue $ws.Range("E13") '  -2.06%  '
Set-TextValue $ws.Range("D14") '3.035.57'
Set-TextValue $ws.Range("E14") '  -0.03%  '
Set-TextValue $ws.Range("D15") '62.821.02'
Set-TextValue $ws.Range("E15") '  -0.70%  '
Set-TextValue $ws.Range("D16") '0.0000144'
Set-TextValue $ws.Range("E16") '  -2.39%  '
Set-TextValue $ws.Range("D17") '2.578.62'
Set-TextValue $ws.Range("E17") '  -0.29%  '
Set-TextValue $ws.Range("D18") '11.05'
Set-TextValue $ws.Range("E18") '  -2.84%  '
Set-TextValue $ws.Range("D19") '340.75'
Set-TextValue $ws.Range("E19") '  -0.60%  '
Set-TextValue $ws.Range("E20") '  -2.11%  '
Set-TextValue $ws.Range("E21") '  -3.92%  '
Set-TextValue $ws.Range("E22") '  -0.02%  '
Set-TextValue $ws.Range("D23") '5.72'
Set-TextValue $ws.Range("E23") '  +3.21%  '
Set-TextValue $ws.Range("D24") '67.71'
Set-TextValue $ws.Range("E24") '  +1.23%  '
Set-TextValue $ws.Range("E25") '  +6.69%  '
Set-TextValue $ws.Range("D26") '1.59'
Set-TextValue $ws.Range("E26") '  -3.16%  '
Set-TextValue $ws.Range("E27") '  -3.55%  '
Set-TextValue $ws.Range("D28") '7.99'
Set-TextValue $ws.Range("E28") '  -2.34%  '
Set-TextValue $ws.Range("E29") '  +1.72%  '
Set-TextValue $ws.Range("D30") '8.23'
Set-TextValue $ws.Range("E30") '  -3.29%  '
Set-TextValue $ws.Range("E31") '  -3.28%  '
Set-TextValue $ws.Range("D32") '459.61'
Set-TextValue $ws.Range("E32") '  -1.44%  '
Set-TextValue $ws.Range("D33") '0.0₃0796'
Set-TextValue $ws.Range("E33") '  -3.83%  '
Set-TextValue $ws.Range("E34") '  +1.20%  '
Set-TextValue $ws.Range("D35") '176.69'
Set-TextValue $ws.Range("E35") '  +0.08%  '
Set-TextValue $ws.Range("E36") '  +0.04%  '
Set-TextValue $ws.Range("E37") '  -2.34%  '
Set-TextValue $ws.Range("E38") '  -2.27%  '
Set-TextValue $ws.Range("E39") '  -0.60%  '
Set-TextValue $ws.Range("E40") '  +0.02%  '
Set-TextValue $ws.Range("E41") '  -3.52%  '
Set-TextValue $ws.Range("D42") '39.98'
Set-TextValue $ws.Range("E42") '  +0.82%  '
Set-TextValue $ws.Range("D43") '157.75'
Set-TextValue $ws.Range("E43") '  +3.99%  '
Set-TextValue $ws.Range("E44") '  -3.67%  '
Set-TextValue $ws.Range("D45") '21.18'
Set-TextValue $ws.Range("E45") '  +0.30%  '
Set-TextValue $ws.Range("D46") '0.632'
Set-TextValue $ws.Range("E46") '  +2.59%  '
Set-TextValue $ws.Range("D47") '0.0535'
Set-TextValue $ws.Range("E47") '  -2.89%  '
Set-TextValue $ws.Range("D48") '0.0958'
Set-TextValue $ws.Range("E48") '  -2.40%  '
Set-TextValue $ws.Range("E49") '  -2.16%  '
Set-TextValue $ws.Range("D50") '18.01'
Set-TextValue $ws.Range("E50") '  -2.54%  '
Set-TextValue $ws.Range("D51") '11.40'
Set-TextValue $ws.Range("E51") '  +0.13%  '
